# Set vertical alignment to "bottom" for every cell in the first two
# (header) rows of the table, matching the OOXML diff that inserts
# <w:vAlign w:val="bottom"/> into each of those cells' <w:tcPr>.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$wdCellAlignVerticalBottom = 3

foreach ($rowIndex in 1..2) {
    $row = $table.Rows.Item($rowIndex)
    foreach ($cell in $row.Cells) {
        $cell.VerticalAlignment = $wdCellAlignVerticalBottom
    }
}
